{"js": "// Word JS API (Office.js) script\n// Body is the content of: async (context) => { ... }\n//\n// The document is a daily \"three-digit x one-digit multiplication\" drill\n// sheet: a date heading followed by a table of 25 problems like\n// \"862\u00d79=\". This edit bumps the date by one day and swaps in a fresh\n// batch of 25 problems, each of which is textually unique in the\n// original document, so a simple search-and-replace per old->new pair\n// is unambiguous.\n\nconst replacements = [\n  [\"2024-03-25 Monday\", \"2024-03-26 Tuesday\"],\n  [\"862\u00d79=\", \"514\u00d76=\"],\n  [\"646\u00d79=\", \"846\u00d75=\"],\n  [\"756\u00d76=\", \"169\u00d78=\"],\n  [\"144\u00d73=\", \"123\u00d76=\"],\n  [\"610\u00d79=\", \"202\u00d73=\"],\n  [\"725\u00d78=\", \"544\u00d73=\"],\n  [\"908\u00d73=\", \"975\u00d77=\"],\n  [\"572\u00d79=\", \"411\u00d72=\"],\n  [\"236\u00d76=\", \"474\u00d78=\"],\n  [\"685\u00d76=\", \"930\u00d73=\"],\n  [\"978\u00d73=\", \"939\u00d74=\"],\n  [\"138\u00d72=\", \"114\u00d78=\"],\n  [\"197\u00d74=\", \"144\u00d78=\"],\n  [\"596\u00d78=\", \"570\u00d76=\"],\n  [\"804\u00d76=\", \"325\u00d74=\"],\n  [\"448\u00d76=\", \"788\u00d78=\"],\n  [\"641\u00d76=\", \"895\u00d76=\"],\n  [\"534\u00d76=\", \"130\u00d73=\"],\n  [\"144\u00d76=\", \"647\u00d75=\"],\n  [\"519\u00d74=\", \"340\u00d79=\"],\n  [\"879\u00d76=\", \"944\u00d78=\"],\n  [\"475\u00d72=\", \"250\u00d76=\"],\n  [\"258\u00d72=\", \"757\u00d72=\"],\n  [\"389\u00d77=\", \"556\u00d77=\"],\n  [\"855\u00d74=\", \"167\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop script\n# $word.ActiveDocument is the open document.\n#\n# The document is a daily \"three-digit x one-digit multiplication\" drill\n# sheet: a date heading followed by a table of 25 problems like\n# \"862\u00d79=\". This edit bumps the date by one day and swaps in a fresh\n# batch of 25 problems. Every old value is textually unique in the\n# original document, so Find/Replace per old->new pair is unambiguous.\n#\n# NOTE: this interpreter flattens nested arrays (e.g. `@(@(\"a\",\"b\"))`\n# collapses to a single flat array), so the old/new values are kept in\n# two parallel arrays instead of an array-of-pairs.\n\n$d = $word.ActiveDocument\n\n$olds = @(\n    \"2024-03-25 Monday\",\n    \"862\u00d79=\",\n    \"646\u00d79=\",\n    \"756\u00d76=\",\n    \"144\u00d73=\",\n    \"610\u00d79=\",\n    \"725\u00d78=\",\n    \"908\u00d73=\",\n    \"572\u00d79=\",\n    \"236\u00d76=\",\n    \"685\u00d76=\",\n    \"978\u00d73=\",\n    \"138\u00d72=\",\n    \"197\u00d74=\",\n    \"596\u00d78=\",\n    \"804\u00d76=\",\n    \"448\u00d76=\",\n    \"641\u00d76=\",\n    \"534\u00d76=\",\n    \"144\u00d76=\",\n    \"519\u00d74=\",\n    \"879\u00d76=\",\n    \"475\u00d72=\",\n    \"258\u00d72=\",\n    \"389\u00d77=\",\n    \"855\u00d74=\"\n)\n\n$news = @(\n    \"2024-03-26 Tuesday\",\n    \"514\u00d76=\",\n    \"846\u00d75=\",\n    \"169\u00d78=\",\n    \"123\u00d76=\",\n    \"202\u00d73=\",\n    \"544\u00d73=\",\n    \"975\u00d77=\",\n    \"411\u00d72=\",\n    \"474\u00d78=\",\n    \"930\u00d73=\",\n    \"939\u00d74=\",\n    \"114\u00d78=\",\n    \"144\u00d78=\",\n    \"570\u00d76=\",\n    \"325\u00d74=\",\n    \"788\u00d78=\",\n    \"895\u00d76=\",\n    \"130\u00d73=\",\n    \"647\u00d75=\",\n    \"340\u00d79=\",\n    \"944\u00d78=\",\n    \"250\u00d76=\",\n    \"757\u00d72=\",\n    \"556\u00d77=\",\n    \"167\u00d79=\"\n)\n\nfor ($i = 0; $i -lt $olds.Count; $i++) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $olds[$i]\n    $find.Replacement.Text = $news[$i]\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
